$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the 5th run data (rows 39-43, columns D-G)
$data = @{
    39 = @(1178.72, 1178.72, 1091.44, 1040.1099999999999)
    40 = @(1209.92, 1209.92, 1123.56, 1058.18)
    41 = @(1195.1300000000001, 1195.1300000000001, 1109.03, 1054.07)
    42 = @(1226.33, 1226.33, 1140.1400000000001, 1067.75)
    43 = @(1189.72, 1189.72, 1103.27, 1041.49)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 4).Value = $values[0]
    $ws.Cells.Item($row, 5).Value = $values[1]
    $ws.Cells.Item($row, 6).Value = $values[2]
    $ws.Cells.Item($row, 7).Value = $values[3]
}

# Update the selection to G43
$ws.Range("G43").Select()
